$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45447 -> 45448) for every data row (rows 2 through 29).
$ws.Range("C2:C29").Value = 45448
